$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 817.4737
$v = $ws.Range("H17").Value
Write-Host ("H17 now = {0}" -f $v)
